$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Report")

# Report generated timestamp
$ws.Range("D5").Value = "Report Generated On: 08/16/2025 12:47 AM"

# Foreman initials instead of full name
$ws.Range("G8").Value = "JH"

# Billing period start date corrected
$ws.Range("C10").Value = "07/21/2025 to 07/27/25"

# Job # flagged as no match
$ws.Range("G13").Value = "#NO MATCH-1"
